# "Additional companies sent for questionaire"
# The "Parent company" and "Location County/City" columns are being
# removed from the Known Locomotive List sheet (the parent-company /
# county-city data isn't being collected for the additional companies
# being sent the questionnaire). Deleting the entire columns shifts
# every later column left and Excel keeps the remaining data, column
# widths and data-validation ranges in sync automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B = "Parent company"
$ws.Columns("B").Delete()

# After the delete above, the old "Location County/City" column (E)
# has shifted left to D.
$ws.Columns("D").Delete()

$ws.Range("J4").Select() | Out-Null
